$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are stored as text so numeric-looking
# strings like "60.80" or "0.0960" keep their exact formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '41.848.22'
$ws.Range('D3').Value = '2.231.63'
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '232.75'
$ws.Range('E5').Value = '  +1.59%  '
$ws.Range('D6').Value = '0.623'
$ws.Range('E6').Value = '  -1.43%  '
$ws.Range('D7').Value = '60.80'
$ws.Range('E7').Value = '  -5.78%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').Value = '58.03'
$ws.Range('E10').Value = '  -1.87%  '
$ws.Range('D11').Value = '0.0905'
$ws.Range('E11').Value = '  +4.39%  '
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('D13').Value = '2.562.64'
$ws.Range('E13').Value = '  +0.72%  '
$ws.Range('D14').Value = '15.70'
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('D15').Value = '22.68'
$ws.Range('E15').Value = '  +2.21%  '
$ws.Range('D16').Value = '0.803'
$ws.Range('E16').Value = '  -2.22%  '
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').Value = '2.247.54'
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('D19').Value = '41.734.01'
$ws.Range('E19').Value = '  +2.47%  '
$ws.Range('D20').Value = '0.0₃0909'
$ws.Range('E20').Value = '  +0.77%  '
$ws.Range('D21').Value = '72.56'
$ws.Range('E21').Value = '  -1.81%  '
$ws.Range('D22').Value = '6.12'
$ws.Range('E22').Value = '  -0.60%  '
$ws.Range('D23').Value = '247.77'
$ws.Range('E23').Value = '  -0.89%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('D27').Value = '9.68'
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('D28').Value = '169.25'
$ws.Range('E28').Value = '  -2.12%  '
$ws.Range('D29').Value = '0.142'
$ws.Range('E29').Value = '  -1.83%  '
$ws.Range('D30').Value = '19.92'
$ws.Range('E30').Value = '  -1.99%  '
$ws.Range('E31').Value = '  -2.52%  '
$ws.Range('E32').Value = '  -6.08%  '
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('D34').Value = '5.05'
$ws.Range('E34').Value = '  +5.69%  '
$ws.Range('E35').Value = '  +0.84%  '
$ws.Range('D36').Value = '0.0656'
$ws.Range('E36').Value = '  +3.98%  '
$ws.Range('D37').Value = '6.56'
$ws.Range('E37').Value = '  -8.20%  '
$ws.Range('E38').Value = '  -2.60%  '
$ws.Range('D39').Value = '3.62'
$ws.Range('E39').Value = '  -5.47%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.0241'
$ws.Range('E41').Value = '  +3.30%  '
$ws.Range('B42').Value = 'TerraClassic'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D42').Value = '0.000235'
$ws.Range('E42').Value = '  +13.43%  '
$ws.Range('D43').Value = '8.61'
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('E45').Value = '  -2.65%  '
$ws.Range('D46').Value = '0.0960'
$ws.Range('E46').Value = '  +2.10%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').Value = '4.44'
$ws.Range('E47').Value = '  -9.75%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '1.470.61'
$ws.Range('E48').Value = '  -2.66%  '
$ws.Range('D49').Value = '16.60'
$ws.Range('E49').Value = '  -4.67%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '2.28'
$ws.Range('E50').Value = '  +7.53%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').Value = '1.08'
$ws.Range('E51').Value = '  -2.80%  '
